$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New handoff identifiers (old bcbcd2ce-... UUID -> new b51b8b31-... UUID),
# new xliff hash (20c32e1f... -> 8ffa296a...) and refreshed timestamps, as
# generated by the localization-status CI report for this handoff.
# ---------------------------------------------------------------------------
$oldGuid = "bcbcd2ce-0db8-4db9-ae74-df2c4b3632dd"
$newGuid = "b51b8b31-b051-48e1-8e39-678998f28687"
$newHash = "8ffa296a09d7cd8823cb642e07bd8253609c5060"

$newFileName        = "$newGuid.md"
$newPathAndName      = "e2e\$newGuid.md"
$newZhXlf            = "$newGuid.$newHash.zh-cn.xlf"
$newDeXlf            = "$newGuid.$newHash.de-de.xlf"
$newOverviewDate     = "2016-08-26 19:02:54"
$newZhHandoffDate    = "2016-08-26 19:02:50"
$unsetHandbackDate   = "0001-01-01 00:00:00"

$oldOverviewHyperlinkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bb07d18c45b6c8471ca286cffa8a792fc47f2dfe/e2e/$oldGuid.md"
$oldZhHyperlinkTarget        = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bb07d18c45b6c8471ca286cffa8a792fc47f2dfe/e2e/$oldGuid.md"
$oldDeHyperlinkTarget        = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bb07d18c45b6c8471ca286cffa8a792fc47f2dfe/e2e/$oldGuid.md"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFileName
$wsOverview.Range("G2").Value = $newOverviewDate

# B2 carries the hyperlink whose display text needs to change; rebuild it
# (the simulated Hyperlinks collection deletes per-sheet, so remove + re-add).
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $oldOverviewHyperlinkTarget, "", "", $newPathAndName)

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newFileName
$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = $newZhHandoffDate
$wsZh.Range("K2").Value = $unsetHandbackDate

# Drop the "Latest Target File" hyperlink (I2) entirely and clear I2/J2 -
# there is no longer a handback target file for this handoff.
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $oldZhHyperlinkTarget, "", "", $newFileName)
$wsZh.Range("I2").Style = "Normal"
$wsZh.Range("I2").Value = ""
$wsZh.Range("J2").Value = ""

$wsZh.Columns.Item(9).ColumnWidth = 17.833333333333336
$wsZh.Columns.Item(10).ColumnWidth = 20.833333333333336

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newFileName
$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("H2").Value = $newOverviewDate
$wsDe.Range("K2").Value = $unsetHandbackDate

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $oldDeHyperlinkTarget, "", "", $newFileName)
$wsDe.Range("I2").Style = "Normal"
$wsDe.Range("I2").Value = ""
$wsDe.Range("J2").Value = ""

$wsDe.Columns.Item(9).ColumnWidth = 17.833333333333336
$wsDe.Columns.Item(10).ColumnWidth = 20.833333333333336
